$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mySheet (4)")

# Boosted the enemy fire-rate baseline (Boss01HP col "O", row 2) from 1 -> 2.
# Every other O3:O100 cell is a shared formula that depends on O2/O101, so
# Excel's recalculation engine updates them automatically once O2 changes.
$ws.Range("O2").Value = 2

# Column Q ("EnemyRateOfFireMax" header column, spreadsheet col 17) was
# narrowed slightly during the same editing session.
$ws.Columns.Item(17).ColumnWidth = 16.63

# Leave the selection on the cell that was actually edited, matching the
# saved cursor position.
$ws.Range("O2").Select()
